# 🔄 Actualización automática del tracker
# Fill in missing "resultado" (G) and "profit" (H) values for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 34; Resultado = "Acierto"; Profit = 1 },
    @{ Row = 37; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 38; Resultado = "Acierto"; Profit = 1.62 },
    @{ Row = 43; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 45; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 53; Resultado = "Acierto"; Profit = 1.2 },
    @{ Row = 54; Resultado = "Fallo";   Profit = -1 }
)

foreach ($u in $updates) {
    $ws.Range("G$($u.Row)").Value = $u.Resultado
    $ws.Range("H$($u.Row)").Value = $u.Profit
}
